$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6668309
$ws.Range("J43").Value = 1994.5454
$ws.Range("L43").Value = 1994.5454
$ws.Range("N43").Value = -2132.5454

$ws.Range("H96").Value = 3502.9
$ws.Range("I96").Value = 696
$ws.Range("J96").Value = 5799.4546
$ws.Range("K96").Value = 2088
$ws.Range("L96").Value = 17398.3638
$ws.Range("M96").Value = -715
$ws.Range("N96").Value = -20144.3638

$ws.Range("H135").Value = 669.4286
$ws.Range("I135").Value = 669.4286
$ws.Range("K135").Value = 6024.8574
$ws.Range("M135").Value = -3489.8574

$ws.Range("H138").Value = 8122
$ws.Range("I138").Value = 24018.334
$ws.Range("J138").Value = 3102.1052
$ws.Range("K138").Value = 72055.00199999999
$ws.Range("L138").Value = 9306.3156
$ws.Range("M138").Value = -66915.00199999999
$ws.Range("N138").Value = -19586.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 144600.02
$ws.Range("I32").Value = 169814.84
$ws.Range("K32").Value = 169814.84
$ws.Range("M32").Value = -169527.84

$ws.Range("H52").Value = 30000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H74").Value = 349588.72
$ws.Range("I74").Value = 2462.1365
$ws.Range("J74").Value = 1013657
$ws.Range("K74").Value = 2462.1365
$ws.Range("L74").Value = 1013657
$ws.Range("M74").Value = -1588.1365
$ws.Range("N74").Value = -1015405

$ws.Range("H77").Value = 349588.72
$ws.Range("I77").Value = 2462.1365
$ws.Range("J77").Value = 1013657
$ws.Range("K77").Value = 12310.6825
$ws.Range("L77").Value = 5068285
$ws.Range("M77").Value = -7942.682500000001
$ws.Range("N77").Value = -5077021

$ws.Range("H97").Value = 13849.6
$ws.Range("J97").Value = 8136.3335
$ws.Range("L97").Value = 8136.3335
$ws.Range("N97").Value = -9128.333500000001

$ws.Range("H102").Value = 1215.64
$ws.Range("I102").Value = 1234.3914
$ws.Range("K102").Value = 1234.3914
$ws.Range("M102").Value = 387.6086

$ws.Range("H132").Value = 1526.8966
$ws.Range("I132").Value = 956.61224
$ws.Range("K132").Value = 2869.83672
$ws.Range("M132").Value = -339.8367200000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 760.2963
$ws.Range("I20").Value = 846.5789
$ws.Range("J20").Value = 555.375
$ws.Range("K20").Value = 846.5789
$ws.Range("L20").Value = 555.375
$ws.Range("M20").Value = -599.5789
$ws.Range("N20").Value = -1049.375

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H80").Value = 50002124
$ws.Range("J80").Value = 2347.7693
$ws.Range("L80").Value = 2347.7693
$ws.Range("N80").Value = -4343.7693

$ws.Range("H83").Value = 50002124
$ws.Range("J83").Value = 2347.7693
$ws.Range("L83").Value = 11738.8465
$ws.Range("N83").Value = -21722.8465

$ws.Range("H107").Value = 8003.129
$ws.Range("I107").Value = 9172.962
$ws.Range("J107").Value = 1920
$ws.Range("K107").Value = 9172.962
$ws.Range("L107").Value = 1920
$ws.Range("M107").Value = -7252.962
$ws.Range("N107").Value = -5760

$ws.Range("H134").Value = 21952900
$ws.Range("J134").Value = 60002260
$ws.Range("L134").Value = 180006780
$ws.Range("N134").Value = -180011850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 17858668
$ws.Range("I16").Value = 20409620
$ws.Range("K16").Value = 20409620
$ws.Range("M16").Value = -20409333

$ws.Range("H31").Value = 2784.6775
$ws.Range("I31").Value = 3231.2
$ws.Range("K31").Value = 3231.2
$ws.Range("M31").Value = -2936.2

$ws.Range("H34").Value = 2784.6775
$ws.Range("I34").Value = 3231.2
$ws.Range("K34").Value = 3231.2
$ws.Range("M34").Value = -3029.2

$ws.Range("H113").Value = 17858668
$ws.Range("I113").Value = 20409620
$ws.Range("K113").Value = 20409620
$ws.Range("M113").Value = -20407450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3883001.8
$ws.Range("I4").Value = 5500713.5
$ws.Range("J4").Value = 1571984.9
$ws.Range("K4").Value = 16502140.5
$ws.Range("L4").Value = 4715954.699999999
$ws.Range("M4").Value = -16502028.5
$ws.Range("N4").Value = -4716178.699999999

$ws.Range("H7").Value = 202.625
$ws.Range("I7").Value = 233.8
$ws.Range("K7").Value = 701.4000000000001
$ws.Range("M7").Value = -589.4000000000001

$ws.Range("H46").Value = 137500700
$ws.Range("J46").Value = 200000960
$ws.Range("L46").Value = 600002880
$ws.Range("N46").Value = -600003062

$ws.Range("H75").Value = 28575490
$ws.Range("J75").Value = 35719252
$ws.Range("L75").Value = 107157756
$ws.Range("N75").Value = -107159752

$ws.Range("H78").Value = 28575490
$ws.Range("J78").Value = 35719252
$ws.Range("L78").Value = 321473268
$ws.Range("N78").Value = -321483252

$ws.Range("H138").Value = 3900.4211
$ws.Range("J138").Value = 3266.3333
$ws.Range("L138").Value = 9798.999899999999
$ws.Range("N138").Value = -20078.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5292.857
$ws.Range("I122").Value = 5216.8335
$ws.Range("K122").Value = 15650.5005
$ws.Range("M122").Value = -13200.5005

$ws.Range("H132").Value = 1432331.4
$ws.Range("I132").Value = 28510.25
$ws.Range("K132").Value = 85530.75
$ws.Range("M132").Value = -83000.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4914.7896
$ws.Range("I22").Value = 1974.5
$ws.Range("J22").Value = 5698.8667
$ws.Range("K22").Value = 1974.5
$ws.Range("L22").Value = 5698.8667
$ws.Range("M22").Value = -1679.5
$ws.Range("N22").Value = -6288.8667

$ws.Range("H27").Value = 4914.7896
$ws.Range("I27").Value = 1974.5
$ws.Range("J27").Value = 5698.8667
$ws.Range("K27").Value = 1974.5
$ws.Range("L27").Value = 5698.8667
$ws.Range("M27").Value = -1867.5
$ws.Range("N27").Value = -5912.8667

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H132").Value = 3828.923
$ws.Range("I132").Value = 2797.7
$ws.Range("K132").Value = 8393.099999999999
$ws.Range("M132").Value = -5863.099999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3757
$ws.Range("J62").Value = 3757
$ws.Range("L62").Value = 3757
$ws.Range("N62").Value = -5005

$ws.Range("H65").Value = 3757
$ws.Range("J65").Value = 3757
$ws.Range("L65").Value = 18785
$ws.Range("N65").Value = -25025

$ws.Range("H107").Value = 1589186.9
$ws.Range("I107").Value = 1380.3334
$ws.Range("K107").Value = 4141.0002
$ws.Range("M107").Value = -2221.0002

$ws.Range("H132").Value = 2061.7097
$ws.Range("I132").Value = 1684.762
$ws.Range("J132").Value = 2853.3
$ws.Range("K132").Value = 5054.286
$ws.Range("L132").Value = 8559.900000000001
$ws.Range("M132").Value = -2524.286
$ws.Range("N132").Value = -13619.9
